$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "365 Days as an exchange student" description down to row 6
# (next to article-5.html), and give article-3.html a new description.
$existing = $ws.Range("B4").Text
$ws.Range("B6").Value = $existing
$ws.Range("B4").Value = "Music Fair"

# Update the active selection
$ws.Range("B7").Select()
